# Scheduled runner update: refresh Universalis market-price-derived figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 588.83
$ws.Range("I15").Value = 588.83
$ws.Range("K15").Value = 1766.49
$ws.Range("M15").Value = -1597.49
$ws.Range("H28").Value = 878.0526
$ws.Range("I28").Value = 889.4706
$ws.Range("J28").Value = 868.8095
$ws.Range("K28").Value = 889.4706
$ws.Range("L28").Value = 868.8095
$ws.Range("M28").Value = -404.4706
$ws.Range("N28").Value = -1838.8095
$ws.Range("H62").Value = 4271.7856
$ws.Range("I62").Value = 1599.375
$ws.Range("J62").Value = 4900.5884
$ws.Range("K62").Value = 1599.375
$ws.Range("L62").Value = 4900.5884
$ws.Range("M62").Value = -975.375
$ws.Range("N62").Value = -6148.5884
$ws.Range("H65").Value = 4271.7856
$ws.Range("I65").Value = 1599.375
$ws.Range("J65").Value = 4900.5884
$ws.Range("K65").Value = 7996.875
$ws.Range("L65").Value = 24502.942
$ws.Range("M65").Value = -4876.875
$ws.Range("N65").Value = -30742.942
$ws.Range("H98").Value = 3010.9285
$ws.Range("I98").Value = 2974.2964
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 2974.2964
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -1476.2964
$ws.Range("N98").Value = -6996
$ws.Range("H122").Value = 3010.9285
$ws.Range("I122").Value = 2974.2964
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8922.889200000001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6472.889200000001
$ws.Range("N122").Value = -16900
$ws.Range("H129").Value = 857.25
$ws.Range("I129").Value = 725.2857
$ws.Range("J129").Value = 911.58826
$ws.Range("K129").Value = 2175.8571
$ws.Range("L129").Value = 2734.76478
$ws.Range("M129").Value = 2824.1429
$ws.Range("N129").Value = -12734.76478
$ws.Range("H132").Value = 4360.467
$ws.Range("I132").Value = 5348.8696
$ws.Range("J132").Value = 1112.8572
$ws.Range("K132").Value = 16046.6088
$ws.Range("L132").Value = 3338.5716
$ws.Range("M132").Value = -13516.6088
$ws.Range("N132").Value = -8398.571599999999
$ws.Range("H135").Value = 10601336
$ws.Range("I135").Value = 408.6111
$ws.Range("J135").Value = 37860864
$ws.Range("K135").Value = 3677.4999
$ws.Range("L135").Value = 340747776
$ws.Range("M135").Value = -1142.4999
$ws.Range("N135").Value = -340752846
$ws.Range("H137").Value = 18519696
$ws.Range("I137").Value = 1161.2258
$ws.Range("J137").Value = 43479460
$ws.Range("K137").Value = 3483.6774
$ws.Range("L137").Value = 130438380
$ws.Range("M137").Value = -933.6773999999996
$ws.Range("N137").Value = -130443480
$ws.Range("H138").Value = 2584.758
$ws.Range("I138").Value = 2112.2727
$ws.Range("J138").Value = 3122.4138
$ws.Range("K138").Value = 6336.8181
$ws.Range("L138").Value = 9367.241399999999
$ws.Range("M138").Value = -1196.8181
$ws.Range("N138").Value = -19647.2414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1248.8387
$ws.Range("I122").Value = 1227.3636
$ws.Range("K122").Value = 3682.0908
$ws.Range("M122").Value = -1232.0908
$ws.Range("H132").Value = 5471106
$ws.Range("I132").Value = 6236621.5
$ws.Range("K132").Value = 18709864.5
$ws.Range("M132").Value = -18707334.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2885099.2
$ws.Range("I31").Value = 1006.9375
$ws.Range("J31").Value = 11537376
$ws.Range("K31").Value = 1006.9375
$ws.Range("L31").Value = 11537376
$ws.Range("M31").Value = -711.9375
$ws.Range("N31").Value = -11537966
$ws.Range("H34").Value = 2885099.2
$ws.Range("I34").Value = 1006.9375
$ws.Range("J34").Value = 11537376
$ws.Range("K34").Value = 1006.9375
$ws.Range("L34").Value = 11537376
$ws.Range("M34").Value = -804.9375
$ws.Range("N34").Value = -11537780
$ws.Range("H105").Value = 937.3570999999999
$ws.Range("I105").Value = 845.55554
$ws.Range("J105").Value = 1102.6
$ws.Range("K105").Value = 845.55554
$ws.Range("L105").Value = 1102.6
$ws.Range("M105").Value = 901.44446
$ws.Range("N105").Value = -4596.6
$ws.Range("H132").Value = 1767.88
$ws.Range("I132").Value = 1789.05
$ws.Range("J132").Value = 1683.2
$ws.Range("K132").Value = 5367.15
$ws.Range("L132").Value = 5049.6
$ws.Range("M132").Value = -2837.15
$ws.Range("N132").Value = -10109.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3091607.5
$ws.Range("I103").Value = 5667160.5
$ws.Range("J103").Value = 943.8
$ws.Range("K103").Value = 17001481.5
$ws.Range("L103").Value = 2831.4
$ws.Range("M103").Value = -17000602.5
$ws.Range("N103").Value = -4589.4
$ws.Range("H132").Value = 90910776
$ws.Range("I132").Value = 200000820
$ws.Range("J132").Value = 2411.6667
$ws.Range("K132").Value = 1800007380
$ws.Range("L132").Value = 21705.0003
$ws.Range("M132").Value = -1800004850
$ws.Range("N132").Value = -26765.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8523.571
$ws.Range("I122").Value = 8299.071
$ws.Range("K122").Value = 24897.213
$ws.Range("M122").Value = -22447.213
$ws.Range("H132").Value = 16950702
$ws.Range("I132").Value = 21278070
$ws.Range("J132").Value = 1844
$ws.Range("K132").Value = 63834210
$ws.Range("L132").Value = 5532
$ws.Range("M132").Value = -63831680
$ws.Range("N132").Value = -10592

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3104.647
$ws.Range("I132").Value = 3150.74
$ws.Range("K132").Value = 9452.219999999999
$ws.Range("M132").Value = -6922.219999999999
$ws.Range("H136").Value = 1179.3871
$ws.Range("I136").Value = 624
$ws.Range("J136").Value = 3083.5715
$ws.Range("K136").Value = 1872
$ws.Range("L136").Value = 9250.7145
$ws.Range("M136").Value = 678
$ws.Range("N136").Value = -14350.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9808426
$ws.Range("I132").Value = 10324636
$ws.Range("J132").Value = 425
$ws.Range("K132").Value = 30973908
$ws.Range("L132").Value = 1275
$ws.Range("M132").Value = -30971378
$ws.Range("N132").Value = -6335
$ws.Range("H136").Value = 4615634.5
$ws.Range("I136").Value = 7584.4
$ws.Range("J136").Value = 142857140
$ws.Range("K136").Value = 22753.2
$ws.Range("L136").Value = 428571420
$ws.Range("M136").Value = -20203.2
$ws.Range("N136").Value = -428576520
